# added 4wk low sales check
# Update the "Forecast Comparison" sheet with recalculated Inventory
# Coverage / Seasonality Index values (and a few revised MyForecast
# figures), then refresh the dependent totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet -------------------------------------------

# Row 2 (W10)
$wsForecast.Range("L2").Value = 0.9

# Row 3 (W11)
$wsForecast.Range("D3").Value = 6
$wsForecast.Range("H3").Value = 20.17
$wsForecast.Range("L3").Value = 0.82

# Row 4 (W12)
$wsForecast.Range("D4").Value = 6
$wsForecast.Range("H4").Value = 19.17
$wsForecast.Range("L4").Value = 0.84

# Row 5 (W13)
$wsForecast.Range("H5").Value = 18.17
$wsForecast.Range("L5").Value = 1.15

# Row 6 (W14)
$wsForecast.Range("H6").Value = 17.17
$wsForecast.Range("L6").Value = 1.05

# Row 7 (W15)
$wsForecast.Range("H7").Value = 16.17
$wsForecast.Range("L7").Value = 1.04

# Row 8 (W16)
$wsForecast.Range("H8").Value = 13
$wsForecast.Range("L8").Value = 1.05

# Row 9 (W17)
$wsForecast.Range("D9").Value = 7
$wsForecast.Range("H9").Value = 12
$wsForecast.Range("L9").Value = 1.08

# Row 10 (W18)
$wsForecast.Range("D10").Value = 7
$wsForecast.Range("H10").Value = 11
$wsForecast.Range("L10").Value = 0.95

# Row 11 (W19)
$wsForecast.Range("D11").Value = 7
$wsForecast.Range("H11").Value = 10
$wsForecast.Range("L11").Value = 1.17

# Row 12 (W20)
$wsForecast.Range("D12").Value = 7
$wsForecast.Range("H12").Value = 9
$wsForecast.Range("L12").Value = 1.03

# Row 13 (W21)
$wsForecast.Range("D13").Value = 8
$wsForecast.Range("H13").Value = 7
$wsForecast.Range("L13").Value = 1.11

# Row 14 (W22)
$wsForecast.Range("D14").Value = 8
$wsForecast.Range("H14").Value = 6
$wsForecast.Range("L14").Value = 1.05

# Row 15 (W23)
$wsForecast.Range("D15").Value = 8
$wsForecast.Range("H15").Value = 5
$wsForecast.Range("L15").Value = 0.89

# Row 16 (W24)
$wsForecast.Range("D16").Value = 8
$wsForecast.Range("H16").Value = 4
$wsForecast.Range("L16").Value = 0.96

# Row 17 (W25)
$wsForecast.Range("D17").Value = 8
$wsForecast.Range("H17").Value = 3
$wsForecast.Range("L17").Value = 1.15

# --- Summary sheet ---------------------------------------------------------

$wsSummary.Range("B9").Value  = "111"
$wsSummary.Range("B10").Value = "50"
$wsSummary.Range("B11").Value = "24"
$wsSummary.Range("B14").Value = "6"
